$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.734.69"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "2.521.69"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "316.47"
$ws.Range("E5").Value = "  +4.59%  "

$ws.Range("D6").Value = "96.09"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -0.87%  "

$ws.Range("D10").Value = "36.34"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "7.76"
$ws.Range("E12").Value = "  +3.21%  "

$ws.Range("E13").Value = "  -2.42%  "

$ws.Range("D14").Value = "2.909.63"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Value = "15.53"
$ws.Range("E15").Value = "  +5.49%  "

$ws.Range("D16").Value = "2.498.31"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "0.862"
$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").Value = "42.716.62"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  -2.58%  "

$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").Value = "71.64"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "253.82"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "2.99"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").Value = "27.18"
$ws.Range("E26").Value = "  -1.22%  "

$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").Value = "2.36"
$ws.Range("E28").Value = "  +12.60%  "

$ws.Range("D29").Value = "10.19"
$ws.Range("E29").Value = "  +1.65%  "

$ws.Range("D30").Value = "38.07"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").Value = "5.94"
$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("D32").Value = "156.47"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("D33").Value = "19.49"
$ws.Range("E33").Value = "  +5.59%  "

$ws.Range("D34").Value = "3.33"
$ws.Range("E34").Value = "  +1.05%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "2.10"
$ws.Range("E35").Value = "  -3.57%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0790"
$ws.Range("E36").Value = "  -1.57%  "

$ws.Range("D37").Value = "2.62"
$ws.Range("E37").Value = "  -4.62%  "

$ws.Range("D38").Value = "0.114"
$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").Value = "24.21"
$ws.Range("E40").Value = "  -7.34%  "

$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "2.04"
$ws.Range("E43").Value = "  -2.35%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0305"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").Value = "2.028.99"

$ws.Range("D47").Value = "84.84"
$ws.Range("E47").Value = "  -3.41%  "

$ws.Range("E48").Value = "  -2.09%  "

$ws.Range("D49").Value = "74.77"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("D50").Value = "2.764.14"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").Value = "0.191"
$ws.Range("E51").Value = "  +1.05%  "
